# Update TPM-derived metrics for Il12a-Il12rb2 LR pairs (rows 2-19)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.1001513154359418
$ws.Range("J2").Value = 0.117124939437743
$ws.Range("M2").Value = 0.550501
$ws.Range("N2").Value = 1.101002
$ws.Range("O2").Value = 0.2113997982137229
$ws.Range("P2").Value = 0.1654784354636026
$ws.Range("Q2").Value = 0.102427867563
$ws.Range("R2").Value = 0.614567205378
$ws.Range("S2").Value = 0.02117196787399701
$ws.Range("T2").Value = 0.01938165173192692

# Row 3
$ws.Range("I3").Value = 0.1001513154359418
$ws.Range("J3").Value = 0.117124939437743
$ws.Range("O3").Value = 0.2079833583923455
$ws.Range("P3").Value = 0.2442061986816758
$ws.Range("S3").Value = 0.02082980693177832
$ws.Range("T3").Value = 0.02860263623091273

# Row 4
$ws.Range("I4").Value = 0.1001513154359418
$ws.Range("J4").Value = 0.117124939437743
$ws.Range("M4").Value = 0.327139
$ws.Range("N4").Value = 0.981417
$ws.Range("O4").Value = 0.1256257819474244
$ws.Range("P4").Value = 0.147505045129239
$ws.Range("Q4").Value = 0.060868463757
$ws.Range("R4").Value = 0.5478161738130001
$ws.Range("S4").Value = 0.01258158731470334
$ws.Range("T4").Value = 0.01727651947752368

# Row 5
$ws.Range("I5").Value = 0.1001513154359418
$ws.Range("J5").Value = 0.117124939437743
$ws.Range("M5").Value = 0.6082780000000001
$ws.Range("N5").Value = 1.216556
$ws.Range("O5").Value = 0.2335869443613126
$ws.Range("P5").Value = 0.1828459744249861
$ws.Range("Q5").Value = 0.113178029514
$ws.Range("R5").Value = 0.6790681770840001
$ws.Range("S5").Value = 0.0233940397464476
$ws.Range("T5").Value = 0.02141582368096161

# Row 6
$ws.Range("I6").Value = 0.1001513154359418
$ws.Range("J6").Value = 0.117124939437743
$ws.Range("M6").Value = 0.4011773333333333
$ws.Range("N6").Value = 1.203532
$ws.Range("O6").Value = 0.1540574991046086
$ws.Range("P6").Value = 0.1808884928368709
$ws.Range("Q6").Value = 0.074644258172
$ws.Range("R6").Value = 0.6717983235480001
$ws.Range("S6").Value = 0.01542906118809797
$ws.Range("T6").Value = 0.02118655376850312

# Row 7
$ws.Range("I7").Value = 0.1001513154359418
$ws.Range("J7").Value = 0.117124939437743
$ws.Range("M7").Value = 0.1753756666666667
$ws.Range("N7").Value = 0.526127
$ws.Range("O7").Value = 0.0673466179805858
$ws.Range("P7").Value = 0.07907585346362569
$ws.Range("Q7").Value = 0.032630922667
$ws.Range("R7").Value = 0.293678304003
$ws.Range("S7").Value = 0.006744852380917517
$ws.Range("T7").Value = 0.009261754547915004

# Row 8
$ws.Range("G8").Value = 0.8076995
$ws.Range("H8").Value = 1.615399
$ws.Range("I8").Value = 0.4347568694579387
$ws.Range("J8").Value = 0.338959581867057
$ws.Range("M8").Value = 0.550501
$ws.Range("N8").Value = 1.101002
$ws.Range("O8").Value = 0.2113997982137229
$ws.Range("P8").Value = 0.1654784354636026
$ws.Range("Q8").Value = 0.4446393824495
$ws.Range("R8").Value = 1.778557529798
$ws.Range("S8").Value = 0.09190751447543812
$ws.Range("T8").Value = 0.0560905012927575

# Row 9
$ws.Range("G9").Value = 0.8076995
$ws.Range("H9").Value = 1.615399
$ws.Range("I9").Value = 0.4347568694579387
$ws.Range("J9").Value = 0.338959581867057
$ws.Range("O9").Value = 0.2079833583923455
$ws.Range("P9").Value = 0.2442061986816758
$ws.Range("Q9").Value = 0.4374535492311667
$ws.Range("R9").Value = 2.624721295387
$ws.Range("S9").Value = 0.09042219379400464
$ws.Range("T9").Value = 0.08277603099448429

# Row 10
$ws.Range("G10").Value = 0.8076995
$ws.Range("H10").Value = 1.615399
$ws.Range("I10").Value = 0.4347568694579387
$ws.Range("J10").Value = 0.338959581867057
$ws.Range("M10").Value = 0.327139
$ws.Range("N10").Value = 0.981417
$ws.Range("O10").Value = 0.1256257819474244
$ws.Range("P10").Value = 0.147505045129239
$ws.Range("Q10").Value = 0.2642300067305
$ws.Range("R10").Value = 1.585380040383
$ws.Range("S10").Value = 0.05461667168266787
$ws.Range("T10").Value = 0.04999824842028824

# Row 11
$ws.Range("G11").Value = 0.8076995
$ws.Range("H11").Value = 1.615399
$ws.Range("I11").Value = 0.4347568694579387
$ws.Range("J11").Value = 0.338959581867057
$ws.Range("M11").Value = 0.6082780000000001
$ws.Range("N11").Value = 1.216556
$ws.Range("O11").Value = 0.2335869443613126
$ws.Range("P11").Value = 0.1828459744249861
$ws.Range("Q11").Value = 0.4913058364610001
$ws.Range("R11").Value = 1.965223345844
$ws.Range("S11").Value = 0.10155352867677
$ws.Range("T11").Value = 0.06197739503716788

# Row 12
$ws.Range("G12").Value = 0.8076995
$ws.Range("H12").Value = 1.615399
$ws.Range("I12").Value = 0.4347568694579387
$ws.Range("J12").Value = 0.338959581867057
$ws.Range("M12").Value = 0.4011773333333333
$ws.Range("N12").Value = 1.203532
$ws.Range("O12").Value = 0.1540574991046086
$ws.Range("P12").Value = 0.1808884928368709
$ws.Range("Q12").Value = 0.3240307315446667
$ws.Range("R12").Value = 1.944184389268
$ws.Range("S12").Value = 0.0669775560272388
$ws.Range("T12").Value = 0.06131388789654789

# Row 13
$ws.Range("G13").Value = 0.8076995
$ws.Range("H13").Value = 1.615399
$ws.Range("I13").Value = 0.4347568694579387
$ws.Range("J13").Value = 0.338959581867057
$ws.Range("M13").Value = 0.1753756666666667
$ws.Range("N13").Value = 0.526127
$ws.Range("O13").Value = 0.0673466179805858
$ws.Range("P13").Value = 0.07907585346362569
$ws.Range("Q13").Value = 0.1416508382788333
$ws.Range("R13").Value = 0.8499050296730001
$ws.Range("S13").Value = 0.02927940480181921
$ws.Range("T13").Value = 0.02680351822581124

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8640563333333334
$ws.Range("H14").Value = 2.592169
$ws.Range("I14").Value = 0.4650918151061195
$ws.Range("J14").Value = 0.5439154786951998
$ws.Range("M14").Value = 0.550501
$ws.Range("N14").Value = 1.101002
$ws.Range("O14").Value = 0.2113997982137229
$ws.Range("P14").Value = 0.1654784354636026
$ws.Range("Q14").Value = 0.4756638755563334
$ws.Range("R14").Value = 2.853983253338
$ws.Range("S14").Value = 0.0983203158642878
$ws.Range("T14").Value = 0.09000628243891812

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8640563333333334
$ws.Range("H15").Value = 2.592169
$ws.Range("I15").Value = 0.4650918151061195
$ws.Range("J15").Value = 0.5439154786951998
$ws.Range("O15").Value = 0.2079833583923455
$ws.Range("P15").Value = 0.2442061986816758
$ws.Range("Q15").Value = 0.4679766543774445
$ws.Range("R15").Value = 4.211789889397
$ws.Range("S15").Value = 0.09673135766656256
$ws.Range("T15").Value = 0.1328275314562788

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8640563333333334
$ws.Range("H16").Value = 2.592169
$ws.Range("I16").Value = 0.4650918151061195
$ws.Range("J16").Value = 0.5439154786951998
$ws.Range("M16").Value = 0.327139
$ws.Range("N16").Value = 0.981417
$ws.Range("O16").Value = 0.1256257819474244
$ws.Range("P16").Value = 0.147505045129239
$ws.Range("Q16").Value = 0.2826665248303334
$ws.Range("R16").Value = 2.543998723473
$ws.Range("S16").Value = 0.05842752295005321
$ws.Range("T16").Value = 0.0802302772314271

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8640563333333334
$ws.Range("H17").Value = 2.592169
$ws.Range("I17").Value = 0.4650918151061195
$ws.Range("J17").Value = 0.5439154786951998
$ws.Range("M17").Value = 0.6082780000000001
$ws.Range("N17").Value = 1.216556
$ws.Range("O17").Value = 0.2335869443613126
$ws.Range("P17").Value = 0.1828459744249861
$ws.Range("Q17").Value = 0.5255864583273335
$ws.Range("R17").Value = 3.153518749964001
$ws.Range("S17").Value = 0.108639375938095
$ws.Range("T17").Value = 0.09945275570685656

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.8640563333333334
$ws.Range("H18").Value = 2.592169
$ws.Range("I18").Value = 0.4650918151061195
$ws.Range("J18").Value = 0.5439154786951998
$ws.Range("M18").Value = 0.4011773333333333
$ws.Range("N18").Value = 1.203532
$ws.Range("O18").Value = 0.1540574991046086
$ws.Range("P18").Value = 0.1808884928368709
$ws.Range("Q18").Value = 0.3466398156564445
$ws.Range("R18").Value = 3.119758340908001
$ws.Range("S18").Value = 0.07165088188927178
$ws.Range("T18").Value = 0.09838805117181985

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.8640563333333334
$ws.Range("H19").Value = 2.592169
$ws.Range("I19").Value = 0.4650918151061195
$ws.Range("J19").Value = 0.5439154786951998
$ws.Range("M19").Value = 0.1753756666666667
$ws.Range("N19").Value = 0.526127
$ws.Range("O19").Value = 0.0673466179805858
$ws.Range("P19").Value = 0.07907585346362569
$ws.Range("Q19").Value = 0.1515344554958889
$ws.Range("R19").Value = 1.363810099463
$ws.Range("S19").Value = 0.03132236079784908
$ws.Range("T19").Value = 0.04301058068989945
